# Updated Global Glider Cal and Ingest sheets:
#  - Changed Cal scattering angle (CC_scattering_angle) value to 140
#  - Changed angular resolution (CC_angular_resolution) value to 1.13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Row 2: CC_scattering_angle value (column F) 117 -> 140
$ws.Range("F2").Value = 140

# Row 4: CC_angular_resolution value (column F) 1.08 -> 1.13
$ws.Range("F4").Value = 1.13

# Reflect the selected cell on this sheet as captured in the saved file
$ws.Range("F15").Select()
